# L&T team for disputes
#
# Original run text: "raised with teaching staff: awarded by peer evaluation"
# Target: split into three runs (same rPr: sz=20, szCs=20):
#   "raised with " | "teaching & learning team" | ": awarded by peer evaluation"
#
# We locate just the "teaching staff" substring, force a run split around it
# by nudging a character-formatting property (which makes the interop layer
# carve the hit out into its own run rather than folding the replacement
# back into the surrounding, identically-formatted run), replace its text,
# then restore the formatting to match its neighbours so the three runs end
# up with identical rPr, exactly like the diff.

$d = $word.ActiveDocument
$rng = $d.Content

$found = $rng.Find.Execute("teaching staff", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text 'teaching staff'"
}

# Force the hit to become its own run (distinct formatting prevents the
# interop layer from re-merging it with the unchanged neighbouring text).
$rng.Font.Bold = 1

# Replace the words themselves.
$rng.Text = "teaching & learning team"

# Restore formatting so the new middle run matches the untouched runs
# around it (sz=20 / szCs=20, no bold) -- matching the target diff's rPr.
$rng.Font.Bold = 0
